$d = $word.ActiveDocument

# Replace the docassemble ${ var } template syntax with Jinja-style {{ var }}
# for sued_when_served and answer_deadline in the "Important 30-day deadline"
# paragraph.

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("`${ sued_when_served }", $true, $false, $false, $false, $false, $true, 1, $false, "{{ sued_when_served }}", 2)

$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$find2.Execute("`${ answer_deadline }", $true, $false, $false, $false, $false, $true, 1, $false, "{{ answer_deadline }}", 2)
